$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.934.83'
$ws.Range("E2").Value = '  +2.30%  '
$ws.Range("D3").Value = '3.035.63'
$ws.Range("E3").Value = '  +1.32%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.76'
$ws.Range("E5").Value = '  -0.22%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '153.90'
$ws.Range("E6").Value = '  +6.92%  '
$ws.Range("D8").Value = '3.030.87'
$ws.Range("E8").Value = '  +1.19%  '
$ws.Range("E9").Value = '  -0.63%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.47'
$ws.Range("E10").Value = '  +9.85%  '
$ws.Range("E11").Value = '  +2.62%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.467'
$ws.Range("E12").Value = '  +1.23%  '
$ws.Range("E13").Value = '  +2.54%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.67'
$ws.Range("E14").Value = '  +4.25%  '
$ws.Range("E15").Value = '  +2.13%  '
$ws.Range("D16").Value = '3.537.91'
$ws.Range("E16").Value = '  +1.33%  '
$ws.Range("E17").Value = '  +1.17%  '
$ws.Range("D18").Value = '62.906.54'
$ws.Range("E18").Value = '  +2.27%  '
$ws.Range("D19").Value = '3.034.19'
$ws.Range("E19").Value = '  +1.46%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '454.22'
$ws.Range("E20").Value = '  +0.20%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.31'
$ws.Range("E21").Value = '  +2.37%  '
$ws.Range("E22").Value = '  +1.59%  '
$ws.Range("E23").Value = '  +1.98%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.32'
$ws.Range("E25").Value = '  +7.91%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.33'
$ws.Range("E26").Value = '  +4.27%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.39'
$ws.Range("E27").Value = '  +3.05%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.45'
$ws.Range("E29").Value = '  +4.08%  '
$ws.Range("E30").Value = '  +1.12%  '
$ws.Range("E31").Value = '  +6.78%  '
$ws.Range("E32").Value = '  +0.10%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.68'
$ws.Range("E33").Value = '  +0.71%  '
$ws.Range("E34").Value = '  +2.36%  '
$ws.Range("D35").Value = '0.0₃0869'
$ws.Range("E35").Value = '  +4.48%  '
$ws.Range("E36").Value = '  +2.25%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.94'
$ws.Range("E37").Value = '  +2.75%  '
$ws.Range("E38").Value = '  +11.19%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.12'
$ws.Range("E39").Value = '  +2.76%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '50.58'
$ws.Range("E40").Value = '  +0.69%  '
$ws.Range("E41").Value = '  +3.83%  '
$ws.Range("E42").Value = '  -1.89%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.306'
$ws.Range("E43").Value = '  +13.79%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '42.03'
$ws.Range("E44").Value = '  +7.71%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '397.88'
$ws.Range("E45").Value = '  +0.72%  '
$ws.Range("E46").Value = '  +1.86%  '
$ws.Range("D47").Value = '2.732.64'
$ws.Range("E47").Value = '  +0.53%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '132.23'
$ws.Range("E48").Value = '  -0.88%  '
$ws.Range("E49").Value = '  -0.01%  '
$ws.Range("E50").Value = '  +4.68%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '24.52'
$ws.Range("E51").Value = '  +4.55%  '
